$d = $word.ActiveDocument

# Locate the paragraph that describes creating TechnicienServiceTest /
# addManager / save repository test (the one that needs to be
# highlighted green, per the commit "TechnicienServiceTest Mockito,
# addManager and save repository").
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*TechnicienServiceTest*addManager*repository*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # wdBrightGreen (4) serialises to OOXML <w:highlight w:val="green"/>.
    # Going through .Font (rather than the Range directly) also reaches
    # the paragraph mark's run properties (w:pPr/w:rPr), matching how the
    # rest of the document's highlighted bullet points are formatted.
    $target.Range.Font.HighlightColorIndex = 4
    Write-Host "Highlighted paragraph: $($target.Range.Text)"
} else {
    Write-Host "Target paragraph not found"
}
